# Applies the "Several rows supplemented" edit to the ChartSettings sheet
# of Parameters_FinalEnergy.xlsx: adds a Date_Change value to the existing
# row 2, and appends three new rows (3-5) cloning row 2's pattern with a
# few tweaked values (ID, Date_Change, Grid_Bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChartSettings")

# --- row 2: only the Date_Change cell (B2) is new -------------------------
$b2 = $ws.Range("B2")
$b2.Value2 = 45478
$b2.NumberFormat = "mm-dd-yy"

# --- helper data for the new rows (3-5) ------------------------------------
$rows = @(
    @{ Row = 3; Id = "FinalEnergy.02"; Date = 45541; GridBottom = 0.15 },
    @{ Row = 4; Id = "FinalEnergy.03"; Date = 45541; GridBottom = 0.25 },
    @{ Row = 5; Id = "FinalEnergy.04"; Date = 45541; GridBottom = 0.25 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A$row").Value2 = $r.Id

    $dateCell = $ws.Range("B$row")
    $dateCell.Value2 = $r.Date
    $dateCell.NumberFormat = "mm-dd-yy"

    $ws.Range("D$row").Value2 = "Final Energy Demand"
    $ws.Range("E$row").Value2 = "Endenergiebedarf"
    $ws.Range("G$row").Value2 = "Annual values by energy carrier"
    $ws.Range("H$row").Value2 = "Jährliche Werte je Energieträger"
    $ws.Range("O$row").Value2 = "Final energy related to reference area [kWh/(m²a)]"
    $ws.Range("P$row").Value2 = "Endenergie pro m² Referenzfläche [kWh/(m²a)]"

    $ws.Range("R$row").Value2 = 20
    $ws.Range("S$row").Value2 = 20

    $pctRange = $ws.Range("T$row" + ":W$row")
    $pctRange.NumberFormat = "0%"
    $ws.Range("T$row").Value2 = 0.12
    $ws.Range("U$row").Value2 = $r.GridBottom
    $ws.Range("V$row").Value2 = 0.1
    $ws.Range("W$row").Value2 = 0.1

    $ws.Range("X$row").Value2 = 0
    $ws.Range("Y$row").Value2 = 400

    $ws.Range("AA$row").Value2 = 5
    $ws.Range("AB$row").Value2 = 50

    $ws.Rows.Item($row).AutoFit() | Out-Null
}

# --- restore selection similar to the authored workbook --------------------
$ws.Range("B2").Select()
